$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.668.99"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.583.26"
$ws.Range("E3").Value = "  -3.30%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.16"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.22"
$ws.Range("E8").Value = "  -5.37%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("D12").Value = "1.808.73"
$ws.Range("E12").Value = "  -3.22%  "
$ws.Range("D13").Value = "1.593.55"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("E14").Value = "  -4.27%  "
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("D16").Value = "27.632.25"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.17"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.64"
$ws.Range("E18").Value = "  -4.67%  "
$ws.Range("E19").Value = "  -3.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  -7.70%  "
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -5.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.48"
$ws.Range("E23").Value = "  -7.08%  "
$ws.Range("E24").Value = "  -6.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.60"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.75"
$ws.Range("E27").Value = "  -3.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.09"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -6.05%  "
$ws.Range("D33").Value = "1.386.57"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E34").Value = "  -5.79%  "
$ws.Range("E35").Value = "  -5.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.969"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.818"
$ws.Range("E40").Value = "  -4.30%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.977"
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("E43").Value = "  -4.55%  "
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.44"
$ws.Range("E45").Value = "  -4.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.21"
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("D47").Value = "1.720.16"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.94"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "0.0₆01000"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("E51").Value = "  -1.21%  "
